$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old D/E/F -> E/F/G),
# carrying the header style of the row along with it.
$ws.Columns.Item(4).Insert()

# New column D header
$ws.Range("D1").Value = "PATH"

# Clear out the placeholder "dfsf" values that used to live in column B
# (rows 2-16), leaving the cells blank.
$ws.Range("B2:B16").ClearContents()
